# Auto-generated script applying numeric corrections to Zalera_Profits sheets
# across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(48, 8).Value = 13754
$ws.Cells.Item(48, 10).Value = 13333
$ws.Cells.Item(48, 12).Value = 39999
$ws.Cells.Item(48, 14).Value = -40583

$ws.Cells.Item(56, 8).Value = 13754
$ws.Cells.Item(56, 10).Value = 13333
$ws.Cells.Item(56, 12).Value = 39999
$ws.Cells.Item(56, 14).Value = -41067

$ws.Cells.Item(64, 8).Value = 10481.2
$ws.Cells.Item(64, 9).Value = 12051.25
$ws.Cells.Item(64, 11).Value = 12051.25
$ws.Cells.Item(64, 13).Value = -11803.25

$ws.Cells.Item(67, 8).Value = 10481.2
$ws.Cells.Item(67, 9).Value = 12051.25
$ws.Cells.Item(67, 11).Value = 12051.25
$ws.Cells.Item(67, 13).Value = -11193.25

$ws.Cells.Item(69, 8).Value = 20000
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 20000
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).Value = 60000
$ws.Cells.Item(69, 14).Value = -61748
$ws.Cells.Item(69, 13).ClearContents()

$ws.Cells.Item(72, 8).Value = 20000
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 20000
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).Value = 180000
$ws.Cells.Item(72, 14).Value = -188736
$ws.Cells.Item(72, 13).ClearContents()

$ws.Cells.Item(86, 8).Value = 3250
$ws.Cells.Item(86, 9).Value = 3000
$ws.Cells.Item(86, 11).Value = 3000
$ws.Cells.Item(86, 13).Value = -1877

$ws.Cells.Item(89, 8).Value = 3250
$ws.Cells.Item(89, 9).Value = 3000
$ws.Cells.Item(89, 11).Value = 15000
$ws.Cells.Item(89, 13).Value = -9384

$ws.Cells.Item(92, 8).Value = 1430.8334
$ws.Cells.Item(92, 9).Value = 1520.6364
$ws.Cells.Item(92, 10).Value = 443
$ws.Cells.Item(92, 11).Value = 1520.6364
$ws.Cells.Item(92, 12).Value = 443
$ws.Cells.Item(92, 13).Value = -272.6364000000001
$ws.Cells.Item(92, 14).Value = -2939

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 21260.389
$ws.Cells.Item(32, 9).Value = 22753.12
$ws.Cells.Item(32, 11).Value = 22753.12
$ws.Cells.Item(32, 13).Value = -22466.12

$ws.Cells.Item(33, 8).Value = 7800
$ws.Cells.Item(33, 9).Value = 7800
$ws.Cells.Item(33, 11).Value = 7800
$ws.Cells.Item(33, 13).Value = -7471

$ws.Cells.Item(61, 8).Value = 6223.517
$ws.Cells.Item(61, 9).Value = 5690.1055
$ws.Cells.Item(61, 10).Value = 7237
$ws.Cells.Item(61, 11).Value = 5690.1055
$ws.Cells.Item(61, 12).Value = 7237
$ws.Cells.Item(61, 13).Value = -5478.1055
$ws.Cells.Item(61, 14).Value = -7661

$ws.Cells.Item(97, 8).Value = 1684184.9
$ws.Cells.Item(97, 9).Value = 2179343
$ws.Cells.Item(97, 11).Value = 2179343
$ws.Cells.Item(97, 13).Value = -2178847

$ws.Cells.Item(102, 8).Value = 4113.9165
$ws.Cells.Item(102, 9).Value = 3668.1936
$ws.Cells.Item(102, 11).Value = 3668.1936
$ws.Cells.Item(102, 13).Value = -2046.1936

$ws.Cells.Item(110, 8).Value = 35715076
$ws.Cells.Item(110, 9).Value = 35715076
$ws.Cells.Item(110, 11).Value = 35715076
$ws.Cells.Item(110, 13).Value = -35713031

$ws.Cells.Item(136, 8).Value = 6223.517
$ws.Cells.Item(136, 9).Value = 5690.1055
$ws.Cells.Item(136, 10).Value = 7237
$ws.Cells.Item(136, 11).Value = 17070.3165
$ws.Cells.Item(136, 12).Value = 21711
$ws.Cells.Item(136, 13).Value = -14520.3165
$ws.Cells.Item(136, 14).Value = -26811

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 1079.3636
$ws.Cells.Item(22, 9).Value = 1177.3
$ws.Cells.Item(22, 10).Value = 100
$ws.Cells.Item(22, 11).Value = 1177.3
$ws.Cells.Item(22, 12).Value = 100
$ws.Cells.Item(22, 13).Value = -1004.3
$ws.Cells.Item(22, 14).Value = -446

$ws.Cells.Item(86, 8).Value = 565746.25
$ws.Cells.Item(86, 9).Value = 4742.5
$ws.Cells.Item(86, 11).Value = 4742.5
$ws.Cells.Item(86, 13).Value = -3619.5

$ws.Cells.Item(89, 8).Value = 565746.25
$ws.Cells.Item(89, 9).Value = 4742.5
$ws.Cells.Item(89, 11).Value = 23712.5
$ws.Cells.Item(89, 13).Value = -18096.5

$ws.Cells.Item(94, 8).Value = 1279.5
$ws.Cells.Item(94, 9).Value = 1127.1786
$ws.Cells.Item(94, 11).Value = 1127.1786
$ws.Cells.Item(94, 13).Value = -676.1786

$ws.Cells.Item(105, 8).Value = 58840204
$ws.Cells.Item(105, 9).Value = 58840204
$ws.Cells.Item(105, 11).Value = 58840204
$ws.Cells.Item(105, 13).Value = -58838457

$ws.Cells.Item(134, 8).Value = 17011.666
$ws.Cells.Item(134, 9).Value = 23855.572
$ws.Cells.Item(134, 10).Value = 11023.25
$ws.Cells.Item(134, 11).Value = 71566.716
$ws.Cells.Item(134, 12).Value = 33069.75
$ws.Cells.Item(134, 13).Value = -69031.716
$ws.Cells.Item(134, 14).Value = -38139.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(57, 8).Value = 39998.5
$ws.Cells.Item(57, 10).Value = 39998.5
$ws.Cells.Item(57, 12).Value = 39998.5
$ws.Cells.Item(57, 14).Value = -41118.5

$ws.Cells.Item(86, 8).Value = 8098.2
$ws.Cells.Item(86, 9).Value = 8109.1113
$ws.Cells.Item(86, 11).Value = 8109.1113
$ws.Cells.Item(86, 13).Value = -6986.1113

$ws.Cells.Item(89, 8).Value = 8098.2
$ws.Cells.Item(89, 9).Value = 8109.1113
$ws.Cells.Item(89, 11).Value = 40545.5565
$ws.Cells.Item(89, 13).Value = -34929.5565

$ws.Cells.Item(94, 8).Value = 6315.222
$ws.Cells.Item(94, 10).Value = 1414
$ws.Cells.Item(94, 12).Value = 1414
$ws.Cells.Item(94, 14).Value = -2316

$ws.Cells.Item(99, 8).Value = 10103.5
$ws.Cells.Item(99, 9).Value = 10103.5
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 10103.5
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -8605.5
$ws.Cells.Item(99, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 10103.5
$ws.Cells.Item(126, 9).Value = 10103.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 30310.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -27840.5
$ws.Cells.Item(126, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 36441.29
$ws.Cells.Item(132, 9).Value = 3501
$ws.Cells.Item(132, 11).Value = 10503
$ws.Cells.Item(132, 13).Value = -7973

$ws.Cells.Item(134, 8).Value = 5542.364
$ws.Cells.Item(134, 9).Value = 5814.4443
$ws.Cells.Item(134, 11).Value = 17443.3329
$ws.Cells.Item(134, 13).Value = -14908.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 1712.5
$ws.Cells.Item(19, 9).Value = 2183.3333
$ws.Cells.Item(19, 10).Value = 300
$ws.Cells.Item(19, 11).Value = 6549.999899999999
$ws.Cells.Item(19, 12).Value = 900
$ws.Cells.Item(19, 13).Value = -6375.999899999999
$ws.Cells.Item(19, 14).Value = -1248

$ws.Cells.Item(47, 8).Value = 125874.25
$ws.Cells.Item(47, 10).Value = 2250
$ws.Cells.Item(47, 12).Value = 6750
$ws.Cells.Item(47, 14).Value = -7612

$ws.Cells.Item(62, 8).Value = 2118.2144
$ws.Cells.Item(62, 9).Value = 1582.9412
$ws.Cells.Item(62, 11).Value = 4748.8236
$ws.Cells.Item(62, 13).Value = -4062.8236

$ws.Cells.Item(65, 8).Value = 2118.2144
$ws.Cells.Item(65, 9).Value = 1582.9412
$ws.Cells.Item(65, 11).Value = 14246.4708
$ws.Cells.Item(65, 13).Value = -10814.4708

$ws.Cells.Item(69, 8).Value = 500
$ws.Cells.Item(69, 9).Value = 500
$ws.Cells.Item(69, 11).Value = 1500
$ws.Cells.Item(69, 13).Value = -689

$ws.Cells.Item(72, 8).Value = 500
$ws.Cells.Item(72, 9).Value = 500
$ws.Cells.Item(72, 11).Value = 4500
$ws.Cells.Item(72, 13).Value = -444

$ws.Cells.Item(74, 8).Value = 11090.5
$ws.Cells.Item(74, 10).Value = 11674.5
$ws.Cells.Item(74, 12).Value = 35023.5
$ws.Cells.Item(74, 14).Value = -37145.5

$ws.Cells.Item(77, 8).Value = 11090.5
$ws.Cells.Item(77, 10).Value = 11674.5
$ws.Cells.Item(77, 12).Value = 105070.5
$ws.Cells.Item(77, 14).Value = -115678.5

$ws.Cells.Item(80, 8).Value = 998.5
$ws.Cells.Item(80, 10).Value = 998.5
$ws.Cells.Item(80, 12).Value = 2995.5
$ws.Cells.Item(80, 14).Value = -4867.5

$ws.Cells.Item(82, 8).Value = 8214
$ws.Cells.Item(82, 10).Value = 8214
$ws.Cells.Item(82, 12).Value = 24642
$ws.Cells.Item(82, 14).Value = -25454

$ws.Cells.Item(83, 8).Value = 998.5
$ws.Cells.Item(83, 10).Value = 998.5
$ws.Cells.Item(83, 12).Value = 8986.5
$ws.Cells.Item(83, 14).Value = -18346.5

$ws.Cells.Item(85, 8).Value = 8214
$ws.Cells.Item(85, 10).Value = 8214
$ws.Cells.Item(85, 12).Value = 24642
$ws.Cells.Item(85, 14).Value = -27450

$ws.Cells.Item(92, 8).Value = 3239.4
$ws.Cells.Item(92, 10).Value = 2174
$ws.Cells.Item(92, 12).Value = 6522
$ws.Cells.Item(92, 14).Value = -9018

$ws.Cells.Item(99, 8).Value = 4499.3335
$ws.Cells.Item(99, 9).Value = 1313.3334
$ws.Cells.Item(99, 11).Value = 3940.0002
$ws.Cells.Item(99, 13).Value = -1694.0002

$ws.Cells.Item(101, 8).Value = 9343
$ws.Cells.Item(101, 10).Value = 9343
$ws.Cells.Item(101, 12).Value = 28029
$ws.Cells.Item(101, 14).Value = -32897

$ws.Cells.Item(131, 10).Value = 4635.9165
$ws.Cells.Item(131, 12).Value = 13907.7495
$ws.Cells.Item(131, 14).Value = -23987.7495

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(35, 8).Value = 7999.6665
$ws.Cells.Item(35, 9).Value = 10000
$ws.Cells.Item(35, 10).Value = 6999.5
$ws.Cells.Item(35, 11).Value = 10000
$ws.Cells.Item(35, 12).Value = 6999.5
$ws.Cells.Item(35, 13).Value = -9702
$ws.Cells.Item(35, 14).Value = -7595.5

$ws.Cells.Item(80, 8).Value = 2524.9666
$ws.Cells.Item(80, 10).Value = 2546.4375
$ws.Cells.Item(80, 12).Value = 2546.4375
$ws.Cells.Item(80, 14).Value = -4542.4375

$ws.Cells.Item(83, 8).Value = 2524.9666
$ws.Cells.Item(83, 10).Value = 2546.4375
$ws.Cells.Item(83, 12).Value = 12732.1875
$ws.Cells.Item(83, 14).Value = -22716.1875

$ws.Cells.Item(102, 8).Value = 4749.75
$ws.Cells.Item(102, 9).Value = 5499.6665
$ws.Cells.Item(102, 11).Value = 5499.6665
$ws.Cells.Item(102, 13).Value = -3877.6665

$ws.Cells.Item(113, 8).Value = 14681.417
$ws.Cells.Item(113, 9).Value = 2567.7
$ws.Cells.Item(113, 10).Value = 75250
$ws.Cells.Item(113, 11).Value = 2567.7
$ws.Cells.Item(113, 12).Value = 75250
$ws.Cells.Item(113, 13).Value = -397.6999999999998
$ws.Cells.Item(113, 14).Value = -79590

$ws.Cells.Item(122, 8).Value = 6562.967
$ws.Cells.Item(122, 9).Value = 6928.905
$ws.Cells.Item(122, 11).Value = 20786.715
$ws.Cells.Item(122, 13).Value = -18336.715

$ws.Cells.Item(132, 8).Value = 11591.167
$ws.Cells.Item(132, 10).Value = 13871
$ws.Cells.Item(132, 12).Value = 41613
$ws.Cells.Item(132, 14).Value = -46673

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2375.6667
$ws.Cells.Item(16, 9).Value = 1622.2727
$ws.Cells.Item(16, 11).Value = 1622.2727
$ws.Cells.Item(16, 13).Value = -1452.2727

$ws.Cells.Item(82, 8).Value = 2897.4546
$ws.Cells.Item(82, 9).Value = 1687.5
$ws.Cells.Item(82, 10).Value = 3588.8572
$ws.Cells.Item(82, 11).Value = 1687.5
$ws.Cells.Item(82, 12).Value = 3588.8572
$ws.Cells.Item(82, 13).Value = -1326.5
$ws.Cells.Item(82, 14).Value = -4310.8572

$ws.Cells.Item(85, 8).Value = 2897.4546
$ws.Cells.Item(85, 9).Value = 1687.5
$ws.Cells.Item(85, 10).Value = 3588.8572
$ws.Cells.Item(85, 11).Value = 1687.5
$ws.Cells.Item(85, 12).Value = 3588.8572
$ws.Cells.Item(85, 13).Value = -439.5
$ws.Cells.Item(85, 14).Value = -6084.8572

$ws.Cells.Item(132, 8).Value = 8195.233
$ws.Cells.Item(132, 9).Value = 7312.0454
$ws.Cells.Item(132, 11).Value = 21936.1362
$ws.Cells.Item(132, 13).Value = -19406.1362

$ws.Cells.Item(136, 8).Value = 3957.7715
$ws.Cells.Item(136, 9).Value = 2610.15
$ws.Cells.Item(136, 10).Value = 5754.6
$ws.Cells.Item(136, 11).Value = 7830.450000000001
$ws.Cells.Item(136, 12).Value = 17263.8
$ws.Cells.Item(136, 13).Value = -5280.450000000001
$ws.Cells.Item(136, 14).Value = -22363.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 39499.75
$ws.Cells.Item(28, 10).Value = 39499.75
$ws.Cells.Item(28, 12).Value = 39499.75
$ws.Cells.Item(28, 14).Value = -40195.75

$ws.Cells.Item(33, 8).Value = 24999.5
$ws.Cells.Item(33, 9).Value = 14000
$ws.Cells.Item(33, 10).Value = 28666
$ws.Cells.Item(33, 11).Value = 14000
$ws.Cells.Item(33, 12).Value = 28666
$ws.Cells.Item(33, 13).Value = -13750
$ws.Cells.Item(33, 14).Value = -29166

$ws.Cells.Item(36, 8).Value = 24999.5
$ws.Cells.Item(36, 9).Value = 14000
$ws.Cells.Item(36, 10).Value = 28666
$ws.Cells.Item(36, 11).Value = 14000
$ws.Cells.Item(36, 12).Value = 28666
$ws.Cells.Item(36, 13).Value = -13750
$ws.Cells.Item(36, 14).Value = -29166

$ws.Cells.Item(43, 8).Value = 47249
$ws.Cells.Item(43, 9).Value = 19997.5
$ws.Cells.Item(43, 10).Value = 74500.5
$ws.Cells.Item(43, 11).Value = 19997.5
$ws.Cells.Item(43, 12).Value = 74500.5
$ws.Cells.Item(43, 13).Value = -19848.5
$ws.Cells.Item(43, 14).Value = -74798.5

$ws.Cells.Item(52, 8).Value = 12050.5
$ws.Cells.Item(52, 9).Value = 6502.2
$ws.Cells.Item(52, 10).Value = 39792
$ws.Cells.Item(52, 11).Value = 6502.2
$ws.Cells.Item(52, 12).Value = 39792
$ws.Cells.Item(52, 13).Value = -6276.2
$ws.Cells.Item(52, 14).Value = -40244

$ws.Cells.Item(58, 8).Value = 9152.5
$ws.Cells.Item(58, 9).Value = 9152.5
$ws.Cells.Item(58, 11).Value = 9152.5
$ws.Cells.Item(58, 13).Value = -8844.5

$ws.Cells.Item(107, 8).Value = 2595.111
$ws.Cells.Item(107, 9).Value = 2049.6365
$ws.Cells.Item(107, 11).Value = 6148.9095
$ws.Cells.Item(107, 13).Value = -4228.9095

$ws.Cells.Item(113, 8).Value = 547.1579
$ws.Cells.Item(113, 9).Value = 426
$ws.Cells.Item(113, 11).Value = 1278
$ws.Cells.Item(113, 13).Value = 892

$ws.Cells.Item(121, 8).Value = 104497.5
$ws.Cells.Item(121, 10).Value = 104497.5
$ws.Cells.Item(121, 12).Value = 104497.5
$ws.Cells.Item(121, 14).Value = -107991.5

$ws.Cells.Item(126, 8).Value = 52211
$ws.Cells.Item(126, 9).Value = 59357.723
$ws.Cells.Item(126, 11).Value = 178073.169
$ws.Cells.Item(126, 13).Value = -175603.169

$ws.Cells.Item(132, 8).Value = 4282.3057
$ws.Cells.Item(132, 9).Value = 1457.4546
$ws.Cells.Item(132, 10).Value = 8721.357
$ws.Cells.Item(132, 11).Value = 4372.3638
$ws.Cells.Item(132, 12).Value = 26164.071
$ws.Cells.Item(132, 13).Value = -1842.3638
$ws.Cells.Item(132, 14).Value = -31224.071
